$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17 (item id 38956)
$ws.Range("H17").Value = 835224.5600000001
$ws.Range("I17").Value = 411.07144
$ws.Range("J17").Value = 1231407.2
$ws.Range("K17").Value = 1233.21432
$ws.Range("L17").Value = 3694221.6
$ws.Range("M17").Value = -1065.21432
$ws.Range("N17").Value = -3694557.6

# Row 106 (item id 19903)
$ws.Range("H106").Value = 4351574
$ws.Range("I106").Value = 5559400.5
$ws.Range("K106").Value = 5559400.5
$ws.Range("M106").Value = -5558769.5

# Row 127 (item id 36114)
$ws.Range("H127").Value = 866.64514
$ws.Range("I127").Value = 349.8
$ws.Range("J127").Value = 1351.1875
$ws.Range("K127").Value = 1049.4
$ws.Range("L127").Value = 4053.5625
$ws.Range("M127").Value = 3910.6
$ws.Range("N127").Value = -13973.5625

# Row 131 (item id 36108)
$ws.Range("H131").Value = 923.625
$ws.Range("I131").Value = 256
$ws.Range("J131").Value = 2036.3334
$ws.Range("K131").Value = 768
$ws.Range("L131").Value = 6109.0002
$ws.Range("M131").Value = 4272
$ws.Range("N131").Value = -16189.0002

# Row 132 (item id 44049)
$ws.Range("H132").Value = 2072.2124
$ws.Range("I132").Value = 2034.3188
$ws.Range("J132").Value = 2309.9092
$ws.Range("K132").Value = 6102.9564
$ws.Range("L132").Value = 6929.7276
$ws.Range("M132").Value = -3572.9564
$ws.Range("N132").Value = -11989.7276


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item id 44147)
$ws.Range("H32").Value = 709736.5
$ws.Range("I32").Value = 839352.7
$ws.Range("J32").Value = 15364.071
$ws.Range("K32").Value = 839352.7
$ws.Range("L32").Value = 15364.071
$ws.Range("M32").Value = -839065.7
$ws.Range("N32").Value = -15938.071

# Row 45 (item id 27714)
$ws.Range("H45").Value = 2394.6
$ws.Range("I45").Value = 2588.8
$ws.Range("J45").Value = 2200.4
$ws.Range("K45").Value = 2588.8
$ws.Range("L45").Value = 2200.4
$ws.Range("M45").Value = -2211.8
$ws.Range("N45").Value = -2954.4

# Row 88 (item id 12530)
$ws.Range("H88").Value = 1886.5333
$ws.Range("I88").Value = 1849.8
$ws.Range("J88").Value = 1960
$ws.Range("K88").Value = 1849.8
$ws.Range("L88").Value = 1960
$ws.Range("M88").Value = -1443.8
$ws.Range("N88").Value = -2772

# Row 91 (item id 12530)
$ws.Range("H91").Value = 1886.5333
$ws.Range("I91").Value = 1849.8
$ws.Range("J91").Value = 1960
$ws.Range("K91").Value = 1849.8
$ws.Range("L91").Value = 1960
$ws.Range("M91").Value = -445.8
$ws.Range("N91").Value = -4768

# Row 132 (item id 43997)
$ws.Range("H132").Value = 2291.525
$ws.Range("I132").Value = 2275.875
$ws.Range("J132").Value = 2328.0417
$ws.Range("K132").Value = 6827.625
$ws.Range("L132").Value = 6984.125100000001
$ws.Range("M132").Value = -4297.625
$ws.Range("N132").Value = -12044.1251


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105 (item id 19947)
$ws.Range("H105").Value = 12501959
$ws.Range("I105").Value = 15626599
$ws.Range("J105").Value = 3400
$ws.Range("K105").Value = 15626599
$ws.Range("L105").Value = 3400
$ws.Range("M105").Value = -15624852
$ws.Range("N105").Value = -6894

# Row 134 (item id 43998)
$ws.Range("H134").Value = 2593.7222
$ws.Range("I134").Value = 2421.6296
$ws.Range("J134").Value = 3110
$ws.Range("K134").Value = 7264.888800000001
$ws.Range("L134").Value = 9330
$ws.Range("M134").Value = -4729.888800000001
$ws.Range("N134").Value = -14400


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (item id 44023)
$ws.Range("H31").Value = 3357.4138
$ws.Range("I31").Value = 1012.451
$ws.Range("J31").Value = 6679.4443
$ws.Range("K31").Value = 1012.451
$ws.Range("L31").Value = 6679.4443
$ws.Range("M31").Value = -717.451
$ws.Range("N31").Value = -7269.4443

# Row 34 (item id 44023)
$ws.Range("H34").Value = 3357.4138
$ws.Range("I34").Value = 1012.451
$ws.Range("J34").Value = 6679.4443
$ws.Range("K34").Value = 1012.451
$ws.Range("L34").Value = 6679.4443
$ws.Range("M34").Value = -810.451
$ws.Range("N34").Value = -7083.4443

# Row 94 (item id 32934)
$ws.Range("H94").Value = 1198.4
$ws.Range("I94").Value = 760
$ws.Range("K94").Value = 760
$ws.Range("M94").Value = -309

# Row 134 (item id 44020)
$ws.Range("H134").Value = 2245.859
$ws.Range("I134").Value = 2254.0806
$ws.Range("J134").Value = 2214
$ws.Range("K134").Value = 6762.2418
$ws.Range("L134").Value = 6642
$ws.Range("M134").Value = -4227.2418
$ws.Range("N134").Value = -11712


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 93 (item id 19808)
$ws.Range("H93").Value = 2873.25
$ws.Range("J93").Value = 4923
$ws.Range("L93").Value = 14769
$ws.Range("N93").Value = -18513

# Row 97 (item id 19846)
$ws.Range("H97").Value = 1225.6666
$ws.Range("J97").Value = 1420.8
$ws.Range("L97").Value = 4262.4
$ws.Range("N97").Value = -5254.4

# Row 138 (item id 44105)
$ws.Range("H138").Value = 5451.769
$ws.Range("J138").Value = 14233.25
$ws.Range("L138").Value = 42699.75
$ws.Range("N138").Value = -52979.75


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 95 (item id 18235)
$ws.Range("H95").Value = 77086
$ws.Range("J95").Value = 77086
$ws.Range("L95").Value = 77086
$ws.Range("N95").Value = -82578

# Row 122 (item id 36182)
$ws.Range("H122").Value = 1453.1333
$ws.Range("I122").Value = 1215.1538
$ws.Range("K122").Value = 3645.4614
$ws.Range("M122").Value = -1195.4614


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16 (item id 5289)
$ws.Range("H16").Value = 1129.0385
$ws.Range("I16").Value = 1181.7368
$ws.Range("K16").Value = 1181.7368
$ws.Range("M16").Value = -1011.7368

# Row 22 (item id 5277)
$ws.Range("H22").Value = 3213.1892
$ws.Range("I22").Value = 396.27777
$ws.Range("J22").Value = 5881.8423
$ws.Range("K22").Value = 396.27777
$ws.Range("L22").Value = 5881.8423
$ws.Range("M22").Value = -101.27777
$ws.Range("N22").Value = -6471.8423

# Row 27 (item id 5277)
$ws.Range("H27").Value = 3213.1892
$ws.Range("I27").Value = 396.27777
$ws.Range("J27").Value = 5881.8423
$ws.Range("K27").Value = 396.27777
$ws.Range("L27").Value = 5881.8423
$ws.Range("M27").Value = -289.27777
$ws.Range("N27").Value = -6095.8423

# Row 32 (item id 2250)
$ws.Range("H32").Value = 1006.5
$ws.Range("I32").Value = 1006.5
$ws.Range("K32").Value = 1006.5
$ws.Range("M32").Value = -689.5

# Row 122 (item id 36247)
$ws.Range("H122").Value = 4579.643
$ws.Range("I122").Value = 3900
$ws.Range("J122").Value = 4957.222
$ws.Range("K122").Value = 11700
$ws.Range("L122").Value = 14871.666
$ws.Range("M122").Value = -9250
$ws.Range("N122").Value = -19771.666


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 53 (item id 3172)
$ws.Range("H53").Value = 100000000
$ws.Range("J53").Value = 100000000
$ws.Range("L53").Value = 100000000
$ws.Range("N53").Value = -100001214

# Row 55 (item id 2832)
$ws.Range("H55").Value = 88766.664
$ws.Range("I55").Value = 3800
$ws.Range("J55").Value = 131250
$ws.Range("K55").Value = 3800
$ws.Range("L55").Value = 131250
$ws.Range("M55").Value = -3523
$ws.Range("N55").Value = -131804

# Row 132 (item id 44029)
$ws.Range("H132").Value = 3942754.5
$ws.Range("I132").Value = 1289.3043
$ws.Range("J132").Value = 10418019
$ws.Range("K132").Value = 3867.9129
$ws.Range("L132").Value = 31254057
$ws.Range("M132").Value = -1337.9129
$ws.Range("N132").Value = -31259117

# Row 136 (item id 44031)
$ws.Range("H136").Value = 654.95
$ws.Range("I136").Value = 568.4342
$ws.Range("K136").Value = 1705.3026
$ws.Range("M136").Value = 844.6974

